$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 5180   # 南昌·CM04动漫游戏博览会
$ws1.Range("F8").Value  = 5309   # 南昌·云芽动漫音乐嘉年华
$ws1.Range("F9").Value  = 614    # 南昌·云芽动漫音乐嘉年华·封茗囧菌内场票
$ws1.Range("F10").Value = 1351   # 南昌·萌卡动漫展
$ws1.Range("F11").Value = 102    # 九江·第二届异次元动漫嘉年华

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 5180   # 南昌·CM04动漫游戏博览会
$ws4.Range("F9").Value  = 5309   # 南昌·云芽动漫音乐嘉年华
$ws4.Range("F10").Value = 614    # 南昌·云芽动漫音乐嘉年华·封茗囧菌内场票
$ws4.Range("F11").Value = 1351   # 南昌·萌卡动漫展
$ws4.Range("F12").Value = 102    # 九江·第二届异次元动漫嘉年华

$wb.Save()
